$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 1127, shifting the existing rows (old 1127-1145) down to 1132-1150
$ws.Rows("1127:1131").Insert()

# Common fixed values shared by the new rows
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria = "Nectarín"
$unidad = "$/bins (420 kilos)"
$kgUnidad = 420

# New row data: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Origen, PrecioKg
$nuevasFilas = @(
    @{ Fila=1127; Fecha=44595; Variedad="June Pearl"; Calidad="Especial"; Volumen=5;  PMin=330000; PMax=330000; PProm=330000; Origen="Región de O'Higgins"; PKg=786 },
    @{ Fila=1128; Fecha=44595; Variedad="June Pearl"; Calidad="Primera";  Volumen=8;  PMin=300000; PMax=300000; PProm=300000; Origen="Región de O'Higgins"; PKg=714 },
    @{ Fila=1129; Fecha=44595; Variedad="June Pearl"; Calidad="Segunda";  Volumen=10; PMin=270000; PMax=270000; PProm=270000; Origen="Región de O'Higgins"; PKg=643 },
    @{ Fila=1130; Fecha=44595; Variedad="Venus";      Calidad="Especial"; Volumen=26; PMin=300000; PMax=330000; PProm=320769; Origen="Región de O'Higgins"; PKg=764 },
    @{ Fila=1131; Fecha=44595; Variedad="Venus";      Calidad="Primera";  Volumen=35; PMin=270000; PMax=290000; PProm=280857; Origen="Región de O'Higgins"; PKg=669 }
)

foreach ($f in $nuevasFilas) {
    $r = $f.Fila
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $f.Fecha
    $ws.Cells.Item($r, 4).Style = $ws.Cells.Item($r + 5, 4).Style
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Cells.Item($r + 5, 4).NumberFormat
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $f.Variedad
    $ws.Cells.Item($r, 12).Value = $f.Calidad
    $ws.Cells.Item($r, 13).Value = $f.Volumen
    $ws.Cells.Item($r, 14).Value = $f.PMin
    $ws.Cells.Item($r, 15).Value = $f.PMax
    $ws.Cells.Item($r, 16).Value = $f.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $f.Origen
    $ws.Cells.Item($r, 19).Value = $f.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
